# WR_89804587_WeekEnding_071325.xlsx edit script
# Applies the changes described by the commit:
#  - refresh the "Report Generated On" timestamp
#  - refresh summary totals (Total Billed Amount, Total Line Items)
#  - clear the Scope ID # value
#  - populate pricing for the Friday / Saturday / Sunday detail rows + day totals
#  - remove the extra "Point 09 / ANC-DSC-16-96-D1" Sunday line item, and shift
#    the remaining "Point 05 / ANC-EXP-8-72-S58" line + TOTAL row up one row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary block -------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:00 AM"

$ws.Range("C8").Value = 11167.04   # Total Billed Amount
$ws.Range("C9").Value = 9          # Total Line Items

$ws.Range("G10").Value = ""        # Scope ID # cleared

# --- Friday (07/11/2025) detail rows ----------------------------------------
$ws.Range("H16").Value = 1297.06
$ws.Range("H17").Value = 1297.06
$ws.Range("H18").Value = 1297.06
$ws.Range("H19").Value = 3891.18   # Friday TOTAL

# --- Saturday (07/12/2025) detail rows --------------------------------------
$ws.Range("H24").Value = 1297.06
$ws.Range("H25").Value = 1297.06
$ws.Range("H26").Value = 1297.06
$ws.Range("H27").Value = 3891.18   # Saturday TOTAL

# --- Sunday (07/13/2025) detail rows ----------------------------------------
$ws.Range("H32").Value = 1297.06
$ws.Range("H33").Value = 1297.06

# Row 34 ("Point 09" / ANC-DSC-16-96-D1) is being dropped from the report.
# Overwrite it in place (preserving its existing row style) with what used to
# be row 35's line item ("Point 05" / ANC-EXP-8-72-S58), updating its pricing.
$ws.Range("A34").Value = "Point 05"
$ws.Range("B34").Value = "ANC-EXP-8-72-S58"
$ws.Range("D34").Value = "ANC,Expanding,8in,72in,Sg Eye 5/8in"
$ws.Range("F34").Value = 1
$ws.Range("H34").Value = 790.5599999999999

# Now remove the old row 35 entirely - this shifts the old row 36 (TOTAL)
# up into row 35, carrying its own style along and fixing up the merged
# cell reference (A36:G36 -> A35:G35) automatically.
$ws.Rows("35").Delete()

# Update the (now-relocated) Sunday TOTAL value.
$ws.Range("H35").Value = 3384.68
